$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update status column C for rows 5 and 6, and add new entry at C8
$ws.Range("C5").Value = "Leido"
$ws.Range("C6").Value = "Leido"
$ws.Range("C8").Value = "enProceso"

# Update active selection to C8
$ws.Range("C8").Select()
